$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy header style from an existing header cell (AC1) to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill the season record values for each data row (2-42)
for ($row = 2; $row -le 42; $row++) {
    $ws.Cells.Item($row, 30).Value = 103   # AD
    $ws.Cells.Item($row, 31).Value = 59    # AE
    $ws.Cells.Item($row, 32).Value = 0     # AF
}
